$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Remove the two mailto hyperlinks (C2, C3)
$ws.Range("C2").Hyperlinks.Delete()
$ws.Range("C3").Hyperlinks.Delete()

# 2. Drop row 3 (the second contact) entirely
$ws.Rows.Item(3).Delete()

# 3. Drop column C (the old "email" column). Deleting column C directly
#    leaves a stale zero-width <col> entry behind in this engine, so we
#    insert a throwaway column ahead of it and delete the pair together,
#    which removes the column definition cleanly.
$ws.Columns.Item(3).Insert()
$ws.Range("C1:D1").EntireColumn.Delete()

# 4. Re-purpose the remaining two columns as "names" / "email".
#    Order matches the authoring order implied by the shared-string table.
$ws.Range("B1").Value = "email"
$ws.Range("B2").Value = "sebastien.debeauffort@outlook.com"
$ws.Range("A1").Value = "names"
$ws.Range("A2").Value = "Jean exemple"

# 5. Column widths (best-fit sized on the new content)
$ws.Columns.Item(1).ColumnWidth = 18.95
$ws.Columns.Item(2).ColumnWidth = 30.8

# 6. B2 keeps the "hyperlink" look even though it's no longer a live link
$ws.Range("B2").Style = "Lien hypertexte"

# 7. Selection moves to A2
$ws.Range("A2").Select() | Out-Null
